$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text representation (avoid Excel auto-converting
# numeric-looking / percentage-looking strings into numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.277.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.181.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.78"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0912"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.77"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.507.69"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.169.44"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.766"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.153.87"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.65"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.88"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.93"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.48%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.35%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.47"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.82%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +12.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.20"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.02"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0807"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.13"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0336"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.97"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "59.53"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.196"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.60"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +12.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.467"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +14.66%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.57%  "
